# Efetuado correções leebank de acordo com o pdf "Apresentação problemas do site lee bank"
#
# Adds a new translation row (id / pt / en / cn) for the "find in a click"
# feature to the Planilha1 sheet, and updates the sheet view / selection
# state to reflect where the author ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- Append the new translation entries (row 44) ---------------------------
# Existing columns are: A = id, B = pt, C = en, D = cn
$ws.Range("A44").Value = "findInAClick"
$ws.Range("B44").Value = "encontre num clique"
$ws.Range("C44").Value = "find in a click"
# cn column keeps reusing the same placeholder string used by the rows above
# (row 38 through 43), which is the shared string "??".
$ws.Range("D44").Value = $ws.Range("D43").Value2

# --- Reflect the resulting view/selection state -----------------------------
$win = $excel.ActiveWindow
$ws.Range("A34").Select() | Out-Null
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("C47").Select() | Out-Null

$wb.Save()
